$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.546.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Value = "'3.329.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'544.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "'171.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.41%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.35%  "

$ws.Range("D8").Value = "'3.322.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.70%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").Value = "'53.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "'8.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.06%  "

$ws.Range("D15").Value = "'3.863.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'17.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.313.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.12%  "

$ws.Range("D19").Value = "'11.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "

$ws.Range("D20").Value = "'63.505.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "

$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").Value = "'408.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").Value = "'4.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.74%  "

$ws.Range("D25").Value = "'13.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.59%  "

$ws.Range("D26").Value = "'82.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("D27").Value = "'10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("D29").Value = "'8.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("D30").Value = "'28.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("D31").Value = "'6.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("D32").Value = "'11.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.45%  "

$ws.Range("D33").Value = "'572.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.97%  "

$ws.Range("E34").Value = "  -2.04%  "

$ws.Range("D35").Value = "'57.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.41%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("E37").Value = "  +1.24%  "

$ws.Range("D38").Value = "'35.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.72%  "

$ws.Range("D39").Value = "'3.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("E40").Value = "  -5.01%  "

$ws.Range("E41").Value = "  -3.05%  "

$ws.Range("D42").Value = "'3.126.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").Value = "'3.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("E46").Value = "  -2.41%  "

$ws.Range("E47").Value = "  -4.64%  "

$ws.Range("E48").Value = "  -4.03%  "

$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("D50").Value = "'131.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("D51").Value = "'8.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "

